$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 28; existing rows 28-37 shift down to 29-38.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new weekly price entry.
$ws.Range("A28").Value = 3
$ws.Range("B28").Value = "Femacal de La Calera"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 44523
$ws.Range("E28").Value = 5
$ws.Range("F28").Value = 100112022
$ws.Range("G28").Value = "Arveja Verde"
$ws.Range("H28").Value = "Perfection"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 70
$ws.Range("K28").Value = 16000
$ws.Range("L28").Value = 16500
$ws.Range("M28").Value = 16250
$ws.Range("N28").Value = "$/malla 25 kilos"
$ws.Range("O28").Value = "Provincia de Talca"
$ws.Range("P28").Value = 650
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"

Write-Output "Row 28 inserted and populated"
